# Auto-generated edit script applying the diff changes to Zalera_Profits workbook
# Updates numeric leve-profit calculation cells across several sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4263
$ws.Range("I40").Value = 3750
$ws.Range("K40").Value = 3750
$ws.Range("M40").Value = -3575
$ws.Range("H87").Value = 60525.58
$ws.Range("J87").Value = 60525.58
$ws.Range("L87").Value = 60525.58
$ws.Range("N87").Value = -63021.58
$ws.Range("H90").Value = 60525.58
$ws.Range("J90").Value = 60525.58
$ws.Range("L90").Value = 181576.74
$ws.Range("N90").Value = -194056.74
$ws.Range("H127").Value = 1051.5
$ws.Range("J127").Value = 5000
$ws.Range("L127").Value = 15000
$ws.Range("N127").Value = -24920
$ws.Range("H132").Value = 828.3103599999999
$ws.Range("I132").Value = 828.1852
$ws.Range("J132").Value = 830
$ws.Range("K132").Value = 2484.5556
$ws.Range("L132").Value = 2490
$ws.Range("M132").Value = 45.44439999999986
$ws.Range("N132").Value = -7550
$ws.Range("H137").Value = 5396.0835
$ws.Range("I137").Value = 1966
$ws.Range("J137").Value = 11112.889
$ws.Range("K137").Value = 5898
$ws.Range("L137").Value = 33338.667
$ws.Range("M137").Value = -3348
$ws.Range("N137").Value = -38438.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32398.875
$ws.Range("I32").Value = 32525.158
$ws.Range("K32").Value = 32525.158
$ws.Range("M32").Value = -32238.158
$ws.Range("H45").Value = 1971.2858
$ws.Range("I45").Value = 1199.75
$ws.Range("K45").Value = 1199.75
$ws.Range("M45").Value = -822.75
$ws.Range("H61").Value = 6583.9
$ws.Range("I61").Value = 5205.9375
$ws.Range("J61").Value = 12095.75
$ws.Range("K61").Value = 5205.9375
$ws.Range("L61").Value = 12095.75
$ws.Range("M61").Value = -4993.9375
$ws.Range("N61").Value = -12519.75
$ws.Range("H74").Value = 372661.28
$ws.Range("I74").Value = 455803.03
$ws.Range("J74").Value = 6837.6
$ws.Range("K74").Value = 455803.03
$ws.Range("L74").Value = 6837.6
$ws.Range("M74").Value = -454929.03
$ws.Range("N74").Value = -8585.6
$ws.Range("H77").Value = 372661.28
$ws.Range("I77").Value = 455803.03
$ws.Range("J77").Value = 6837.6
$ws.Range("K77").Value = 2279015.15
$ws.Range("L77").Value = 34188
$ws.Range("M77").Value = -2274647.15
$ws.Range("N77").Value = -42924
$ws.Range("H122").Value = 3182.2
$ws.Range("I122").Value = 2204.111
$ws.Range("K122").Value = 6612.333
$ws.Range("M122").Value = -4162.333
$ws.Range("H136").Value = 6583.9
$ws.Range("I136").Value = 5205.9375
$ws.Range("J136").Value = 12095.75
$ws.Range("K136").Value = 15617.8125
$ws.Range("L136").Value = 36287.25
$ws.Range("M136").Value = -13067.8125
$ws.Range("N136").Value = -41387.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 198.3
$ws.Range("I7").Value = 148
$ws.Range("K7").Value = 148
$ws.Range("M7").Value = -35
$ws.Range("H58").Value = 5057.32
$ws.Range("I58").Value = 3224.2144
$ws.Range("J58").Value = 7390.364
$ws.Range("K58").Value = 3224.2144
$ws.Range("L58").Value = 7390.364
$ws.Range("M58").Value = -3021.2144
$ws.Range("N58").Value = -7796.364
$ws.Range("H105").Value = 576
$ws.Range("I105").Value = 576
$ws.Range("K105").Value = 576
$ws.Range("M105").Value = 1171
$ws.Range("H136").Value = 5057.32
$ws.Range("I136").Value = 3224.2144
$ws.Range("J136").Value = 7390.364
$ws.Range("K136").Value = 9672.643199999999
$ws.Range("L136").Value = 22171.092
$ws.Range("M136").Value = -7122.643199999999
$ws.Range("N136").Value = -27271.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1925.3889
$ws.Range("H65").Value = 1925.3889
$ws.Range("H92").Value = 1619.8
$ws.Range("I92").Value = 1500
$ws.Range("J92").Value = 1699.6666
$ws.Range("K92").Value = 4500
$ws.Range("L92").Value = 5098.9998
$ws.Range("M92").Value = -3252
$ws.Range("N92").Value = -7594.9998
$ws.Range("H113").Value = 778.7857
$ws.Range("I113").Value = 320.2
$ws.Range("J113").Value = 1033.5555
$ws.Range("K113").Value = 960.5999999999999
$ws.Range("L113").Value = 3100.6665
$ws.Range("M113").Value = 1209.4
$ws.Range("N113").Value = -7440.666499999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3329071.5
$ws.Range("J11").Value = 260700.2
$ws.Range("L11").Value = 260700.2
$ws.Range("N11").Value = -260978.2
$ws.Range("H14").Value = 202025.67
$ws.Range("J14").Value = 2723.6
$ws.Range("L14").Value = 2723.6
$ws.Range("N14").Value = -3059.6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H102").Value = 5322.5
$ws.Range("I102").Value = 5763.3335
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 5763.3335
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -4141.3335
$ws.Range("N102").Value = -7244
$ws.Range("H122").Value = 5238.3477
$ws.Range("I122").Value = 6098.6665
$ws.Range("K122").Value = 18295.9995
$ws.Range("M122").Value = -15845.9995
$ws.Range("H132").Value = 6340.3213
$ws.Range("I132").Value = 4739.65
$ws.Range("K132").Value = 14218.95
$ws.Range("M132").Value = -11688.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3824.4546
$ws.Range("I22").Value = 2535.5334
$ws.Range("K22").Value = 2535.5334
$ws.Range("M22").Value = -2240.5334
$ws.Range("H27").Value = 3824.4546
$ws.Range("I27").Value = 2535.5334
$ws.Range("K27").Value = 2535.5334
$ws.Range("M27").Value = -2428.5334
$ws.Range("H46").Value = 8449.1
$ws.Range("J46").Value = 9221.166999999999
$ws.Range("L46").Value = 9221.166999999999
$ws.Range("N46").Value = -9597.166999999999
$ws.Range("H68").Value = 2109.4
$ws.Range("I68").Value = 1871.2858
$ws.Range("J68").Value = 2665
$ws.Range("K68").Value = 1871.2858
$ws.Range("L68").Value = 2665
$ws.Range("M68").Value = -1122.2858
$ws.Range("N68").Value = -4163
$ws.Range("H71").Value = 2109.4
$ws.Range("I71").Value = 1871.2858
$ws.Range("J71").Value = 2665
$ws.Range("K71").Value = 9356.429
$ws.Range("L71").Value = 13325
$ws.Range("M71").Value = -5612.429
$ws.Range("N71").Value = -20813
$ws.Range("H122").Value = 2848.3076
$ws.Range("I122").Value = 2675.5715
$ws.Range("K122").Value = 8026.7145
$ws.Range("M122").Value = -5576.7145
$ws.Range("H132").Value = 5313.4287
$ws.Range("I132").Value = 3265
$ws.Range("K132").Value = 9795
$ws.Range("M132").Value = -7265
$ws.Range("H136").Value = 4622.55
$ws.Range("I136").Value = 3421.9375
$ws.Range("J136").Value = 9425
$ws.Range("K136").Value = 10265.8125
$ws.Range("L136").Value = 28275
$ws.Range("M136").Value = -7715.8125
$ws.Range("N136").Value = -33375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2090.1052
$ws.Range("I107").Value = 731.7692
$ws.Range("K107").Value = 2195.3076
$ws.Range("M107").Value = -275.3076000000001
$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -59800
$ws.Range("H132").Value = 7822.706
$ws.Range("I132").Value = 6806.185
$ws.Range("K132").Value = 20418.555
$ws.Range("M132").Value = -17888.555

